# Update the "scenario" worksheet's Non-residential / Residential values
# with full-precision figures (Austria parquet data) instead of the
# previously rounded numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario")

$ws.Range("B2").Value = 119.68710703486188
$ws.Range("C2").Value = 319.47331104198275

$ws.Range("B3").Value = 119.16288691782539
$ws.Range("C3").Value = 315.32948722099496

$ws.Range("B4").Value = 119.65008370845369
$ws.Range("C4").Value = 307.08429024595165

$ws.Range("B5").Value = 115.46035912077282
$ws.Range("C5").Value = 295.11620028256016

$ws.Range("B6").Value = 110.56468275012404
$ws.Range("C6").Value = 281.5637429123412

$ws.Range("B7").Value = 106.06154966260674
$ws.Range("C7").Value = 267.9345929015514

$ws.Range("B8").Value = 101.86143875879341
$ws.Range("C8").Value = 253.81099237453978
